$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New header cells (row 1) ---
$ws.Range("F1").Value2 = "SALE DATETIME"
$ws.Range("G1").Value2 = "CUSTOMER COMMENT"
$ws.Range("H1").Value2 = "CUSTOMER RATING"

# --- New data values ---
$ws.Range("F10").Value2 = 44260.46597222222
$ws.Range("F11").Value2 = 44258.598611111112
$ws.Range("G11").Value2 = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum."
$ws.Range("H11").Value2 = 3

Write-Host "values done"
